# Apply the "deliveries" edit:
#  1. Update the date in the heading from 18/02/2021 to 19/02/2021.
#  2. Remove the second paragraph (Brian Binks / address / mobile block) entirely.

$d = $word.ActiveDocument

# 1) Fix the date in the heading.
$d.Content.Find.Execute("Deliveres for 18/02/2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Deliveres for 19/02/2021", 2)

# 2) Remove the whole second paragraph (including its trailing paragraph mark).
#    Paragraph 1 is the heading, paragraph 2 is the Brian Binks block.
$p = $d.Paragraphs(2)
$p.Range.Delete()
